$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("Skill Description") before the existing SFIA Level column.
$ws.Columns("B").Insert()

# Header row
$ws.Range("B1").Value = "Skill Description"

# Map of SkillCode -> full skill name (used to populate new column B for data rows)
$skillNames = @{
    "Autonomy"   = "Autonomy"
    "Influence"  = "Influence"
    "Complexity" = "Complexity"
    "Knowledge"  = "Knowledge"
    "SWDN"       = "Software design"
    "PROG"       = "Programming/software development"
    "SINT"       = "Systems integration and build"
    "TEST"       = "Testing"
    "CFMG"       = "Configuration management"
    "MADE"       = "MADE"
    "REQM"       = "Requirements definition and management"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $code = $ws.Cells.Item($r, 1).Value2
    if ($code -ne $null -and $code -ne "") {
        if ($skillNames.ContainsKey($code)) {
            $ws.Cells.Item($r, 2).Value = $skillNames[$code]
        } else {
            $ws.Cells.Item($r, 2).Value = $code
        }
    }
}
